$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common values shared across the three new rows (subset job outputs for
# Sex/Age/BP only, row subset 50%/25%/10% of dataset ukb51139_subset.csv)
$dataUsed   = "ukb51139_subset.csv"
$colSubset  = "Sex/Age/BP only"
$rowSubset  = "no events"
$target     = "> 140/80"
$normalize  = "zscore"
$impute     = "median"
$covariate  = "none"
$maxCPCA    = 50

# Row 20 - 50% subset
$ws.Range("A20").Value = $dataUsed
$ws.Range("B20").Value = "14006 x 4"
$ws.Range("C20").Value = $colSubset
$ws.Range("D20").Value = $rowSubset
$ws.Range("E20").Value = $target
$ws.Range("F20").Value = $normalize
$ws.Range("G20").Value = $impute
$ws.Range("H20").Value = $covariate
$ws.Range("I20").Value = $maxCPCA
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = "50.2 & 16.6"
$ws.Range("M20").Value = "33.1 & 78.8"
$ws.Range("N20").Value = 3
$ws.Range("O20").Value = 1.18

# Row 21 - 25% subset
$ws.Range("A21").Value = $dataUsed
$ws.Range("B21").Value = "7003 x 4"
$ws.Range("C21").Value = $colSubset
$ws.Range("D21").Value = $rowSubset
$ws.Range("E21").Value = $target
$ws.Range("F21").Value = $normalize
$ws.Range("G21").Value = $impute
$ws.Range("H21").Value = $covariate
$ws.Range("I21").Value = $maxCPCA
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = "60.4 & 28.0"
$ws.Range("M21").Value = "82.1 & 119.2"
$ws.Range("N21").Value = 3
$ws.Range("O21").Value = 0.31

# Row 22 - 10% subset
$ws.Range("A22").Value = $dataUsed
$ws.Range("B22").Value = "2801 x 4"
$ws.Range("C22").Value = $colSubset
$ws.Range("D22").Value = $rowSubset
$ws.Range("E22").Value = $target
$ws.Range("F22").Value = $normalize
$ws.Range("G22").Value = $impute
$ws.Range("H22").Value = $covariate
$ws.Range("I22").Value = $maxCPCA
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = "83.5 & 106.7"
$ws.Range("M22").Value = "53.6 & 132.2"
$ws.Range("N22").Value = 3
$ws.Range("O22").Value = 1.02
